$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.246.67"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").Value = "2.605.05"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'541.18"
$ws.Range("E5").Value = "  +4.12%  "

# Row 6
$ws.Range("D6").Value = "'141.77"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.35%  "

# Row 9
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("E10").Value = "  +2.23%  "

# Row 11
$ws.Range("D11").Value = "'0.336"
$ws.Range("E11").Value = "  +1.32%  "

# Row 12
$ws.Range("E12").Value = "  +1.09%  "

# Row 13
$ws.Range("D13").Value = "3.062.73"
$ws.Range("E13").Value = "  +0.67%  "

# Row 14
$ws.Range("D14").Value = "59.164.44"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15
$ws.Range("D15").Value = "'20.63"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16
$ws.Range("D16").Value = "2.642.63"
$ws.Range("E16").Value = "  +2.03%  "

# Row 17
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  +1.03%  "

# Row 18
$ws.Range("D18").Value = "'341.95"
$ws.Range("E18").Value = "  +0.83%  "

# Row 19
$ws.Range("D19").Value = "'4.38"
$ws.Range("E19").Value = "  +1.42%  "

# Row 20
$ws.Range("D20").Value = "'10.15"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").Value = "'6.39"
$ws.Range("E21").Value = "  -1.59%  "

# Row 22
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").Value = "'67.53"
$ws.Range("E23").Value = "  +1.88%  "

# Row 24
$ws.Range("D24").Value = "'0.409"
$ws.Range("E24").Value = "  +0.96%  "

# Row 25
$ws.Range("E25").Value = "  -0.87%  "

# Row 26
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").Value = "'7.22"
$ws.Range("E27").Value = "  +2.44%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0767"
$ws.Range("E28").Value = "  +6.09%  "

# Row 30
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "  +9.29%  "

# Row 31
$ws.Range("D31").Value = "'5.83"
$ws.Range("E31").Value = "  -2.03%  "

# Row 32
$ws.Range("D32").Value = "'18.75"
$ws.Range("E32").Value = "  -0.11%  "

# Row 33
$ws.Range("D33").Value = "'149.52"
$ws.Range("E33").Value = "  +0.37%  "

# Row 34
$ws.Range("D34").Value = "'4.00"
$ws.Range("E34").Value = "  +0.50%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.12"
$ws.Range("E35").Value = "  -0.90%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'37.17"
$ws.Range("E36").Value = "  +2.42%  "

# Row 37
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  +0.53%  "

# Row 38
$ws.Range("D38").Value = "'0.840"
$ws.Range("E38").Value = "  +1.01%  "

# Row 39
$ws.Range("D39").Value = "'0.824"
$ws.Range("E39").Value = "  +0.62%  "

# Row 40
$ws.Range("D40").Value = "'3.58"
$ws.Range("E40").Value = "  +2.04%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").Value = "'274.61"
$ws.Range("E42").Value = "  -0.41%  "

# Row 43
$ws.Range("E43").Value = "  +1.53%  "

# Row 44
$ws.Range("E44").Value = "  -0.09%  "

# Row 45
$ws.Range("D45").Value = "'0.0958"
$ws.Range("E45").Value = "  +0.74%  "

# Row 46
$ws.Range("D46").Value = "'0.0526"
$ws.Range("E46").Value = "  +0.87%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'18.60"
$ws.Range("E47").Value = "  +3.82%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.951.46"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0224"
$ws.Range("E49").Value = "  +1.69%  "

# Row 50
$ws.Range("D50").Value = "'4.52"
$ws.Range("E50").Value = "  +0.46%  "

# Row 51
$ws.Range("D51").Value = "'111.45"
$ws.Range("E51").Value = "  -1.21%  "
